$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = '64.804.86'
$ws.Cells.Item(2, 5).Value = '  +0.03%  '
$ws.Cells.Item(3, 4).Value = '3.150.29'
$ws.Cells.Item(3, 5).Value = '  -0.20%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Value = '579.03'
$ws.Cells.Item(5, 5).Value = '  +1.23%  '
$ws.Cells.Item(6, 4).Value = '149.25'
$ws.Cells.Item(6, 5).Value = '  -1.07%  '
$ws.Cells.Item(7, 5).Value = '  +0.05%  '
$ws.Cells.Item(8, 4).Value = '3.146.99'
$ws.Cells.Item(8, 5).Value = '  -0.29%  '
$ws.Cells.Item(9, 4).Value = '0.527'
$ws.Cells.Item(10, 4).Value = '0.159'
$ws.Cells.Item(10, 5).Value = '  -2.59%  '
$ws.Cells.Item(11, 4).Value = '6.12'
$ws.Cells.Item(11, 5).Value = '  -1.15%  '
$ws.Cells.Item(12, 4).Value = '0.501'
$ws.Cells.Item(12, 5).Value = '  -1.17%  '
$ws.Cells.Item(13, 4).Value = '0.0000264'
$ws.Cells.Item(13, 5).Value = '  +2.45%  '
$ws.Cells.Item(14, 4).Value = '37.14'
$ws.Cells.Item(14, 5).Value = '  -2.76%  '
$ws.Cells.Item(15, 4).Value = '3.664.50'
$ws.Cells.Item(15, 5).Value = '  -0.21%  '
$ws.Cells.Item(16, 4).Value = '64.816.37'
$ws.Cells.Item(16, 5).Value = '  -0.13%  '
$ws.Cells.Item(17, 2).Value = 'Polkadot'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(17, 4).Value = '7.15'
$ws.Cells.Item(17, 5).Value = '  -1.36%  '
$ws.Cells.Item(18, 2).Value = 'WrappedEther'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(18, 4).Value = '3.144.30'
$ws.Cells.Item(18, 5).Value = '  -0.35%  '
$ws.Cells.Item(19, 5).Value = '  +0.29%  '
$ws.Cells.Item(20, 4).Value = '504.97'
$ws.Cells.Item(20, 5).Value = '  -2.52%  '
$ws.Cells.Item(21, 4).Value = '14.90'
$ws.Cells.Item(21, 5).Value = '  -0.49%  '
$ws.Cells.Item(22, 5).Value = '  -3.21%  '
$ws.Cells.Item(23, 4).Value = '15.17'
$ws.Cells.Item(23, 5).Value = '  -0.93%  '
$ws.Cells.Item(24, 5).Value = '  -2.05%  '
$ws.Cells.Item(25, 4).Value = '84.37'
$ws.Cells.Item(25, 5).Value = '  -1.30%  '
$ws.Cells.Item(26, 5).Value = '  -0.04%  '
$ws.Cells.Item(27, 2).Value = 'RenderToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(27, 4).Value = '8.95'
$ws.Cells.Item(27, 5).Value = '  +1.89%  '
$ws.Cells.Item(28, 2).Value = 'PancakeSwap'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(28, 4).Value = '2.92'
$ws.Cells.Item(28, 5).Value = '  -0.53%  '
$ws.Cells.Item(29, 5).Value = '  -0.49%  '
$ws.Cells.Item(30, 2).Value = 'Stacks'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(30, 4).Value = '2.79'
$ws.Cells.Item(30, 5).Value = '  +3.58%  '
$ws.Cells.Item(31, 2).Value = 'EthereumClassic'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(31, 4).Value = '27.60'
$ws.Cells.Item(31, 5).Value = '  -1.45%  '
$ws.Cells.Item(32, 4).Value = '0.999'
$ws.Cells.Item(32, 5).Value = '  -0.04%  '
$ws.Cells.Item(33, 5).Value = '  +0.54%  '
$ws.Cells.Item(34, 5).Value = '  +1.84%  '
$ws.Cells.Item(35, 4).Value = '6.48'
$ws.Cells.Item(35, 5).Value = '  -1.53%  '
$ws.Cells.Item(36, 4).Value = '54.86'
$ws.Cells.Item(36, 5).Value = '  -1.47%  '
$ws.Cells.Item(37, 4).Value = '484.52'
$ws.Cells.Item(37, 5).Value = '  -1.13%  '
$ws.Cells.Item(38, 5).Value = '  +2.38%  '
$ws.Cells.Item(39, 4).Value = '0.0416'
$ws.Cells.Item(39, 5).Value = '  -1.94%  '
$ws.Cells.Item(40, 5).Value = '  -2.02%  '
$ws.Cells.Item(41, 5).Value = '  +0.64%  '
$ws.Cells.Item(42, 4).Value = '2.987.37'
$ws.Cells.Item(42, 5).Value = '  -4.26%  '
$ws.Cells.Item(43, 5).Value = '  -2.91%  '
$ws.Cells.Item(44, 5).Value = '  -0.98%  '
$ws.Cells.Item(45, 5).Value = '  -4.85%  '
$ws.Cells.Item(46, 4).Value = '28.11'
$ws.Cells.Item(46, 5).Value = '  -4.36%  '
$ws.Cells.Item(47, 4).Value = '0.0₃0588'
$ws.Cells.Item(47, 5).Value = '  +1.56%  '
$ws.Cells.Item(49, 4).Value = '0.114'
$ws.Cells.Item(49, 5).Value = '  -1.85%  '
$ws.Cells.Item(50, 4).Value = '2.24'
$ws.Cells.Item(50, 5).Value = '  -3.10%  '
$ws.Cells.Item(51, 4).Value = '2.47'
$ws.Cells.Item(51, 5).Value = '  +13.62%  '
